$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '61.989.68'
$ws.Range("E2").Value = '  -1.74%  '

# Row 3
$ws.Range("D3").Value = '3.423.68'
$ws.Range("E3").Value = '  -1.01%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '578.68'
$ws.Range("E5").Value = '  -0.31%  '

# Row 6
$ws.Range("D6").Value = '153.22'
$ws.Range("E6").Value = '  +3.88%  '

# Row 7
$ws.Range("E7").Value = '  +0.05%  '

# Row 8
$ws.Range("E8").Value = '  +1.47%  '

# Row 9
$ws.Range("D9").Value = '8.06'
$ws.Range("E9").Value = '  +3.29%  '

# Row 10
$ws.Range("E10").Value = '  +0.28%  '

# Row 11
$ws.Range("E11").Value = '  +3.55%  '

# Row 12
$ws.Range("D12").Value = '4.008.84'
$ws.Range("E12").Value = '  -1.14%  '

# Row 13
$ws.Range("E13").Value = '  +0.76%  '

# Row 14
$ws.Range("D14").Value = '28.51'
$ws.Range("E14").Value = '  -2.74%  '

# Row 15
$ws.Range("D15").Value = '3.420.99'
$ws.Range("E15").Value = '  -1.18%  '

# Row 16
$ws.Range("D16").Value = '0.0000172'
$ws.Range("E16").Value = '  -0.15%  '

# Row 17
$ws.Range("D17").Value = '62.015.74'
$ws.Range("E17").Value = '  -1.74%  '

# Row 18
$ws.Range("D18").Value = '6.56'
$ws.Range("E18").Value = '  +2.28%  '

# Row 19
$ws.Range("E19").Value = '  +0.29%  '

# Row 20
$ws.Range("D20").Value = '8.96'
$ws.Range("E20").Value = '  -3.49%  '

# Row 21
$ws.Range("D21").Value = '382.69'
$ws.Range("E21").Value = '  -1.34%  '

# Row 22
$ws.Range("D22").Value = '0.573'
$ws.Range("E22").Value = '  +1.77%  '

# Row 23
$ws.Range("D23").Value = '75.27'
$ws.Range("E23").Value = '  +1.06%  '

# Row 24
$ws.Range("E24").Value = '  +0.05%  '

# Row 25
$ws.Range("D25").Value = '3.559.65'
$ws.Range("E25").Value = '  -1.36%  '

# Row 26
$ws.Range("E26").Value = '  -2.52%  '

# Row 27
$ws.Range("E27").Value = '  -2.08%  '

# Row 28
$ws.Range("D28").Value = '7.67'
$ws.Range("E28").Value = '  +0.58%  '

# Row 29
$ws.Range("E29").Value = '  +0.04%  '

# Row 30
$ws.Range("E30").Value = '  -0.58%  '

# Row 31
$ws.Range("D31").Value = '7.88'
$ws.Range("E31").Value = '  -3.56%  '

# Row 32
$ws.Range("D32").Value = '0.999'
$ws.Range("E32").Value = '  -0.03%  '

# Row 33
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '23.27'
$ws.Range("E33").Value = '  -0.54%  '

# Row 34
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = '1.33'
$ws.Range("E34").Value = '  -0.70%  '

# Row 35
$ws.Range("D35").Value = '5.53'
$ws.Range("E35").Value = '  +3.74%  '

# Row 36
$ws.Range("E36").Value = '  +0.10%  '

# Row 37
$ws.Range("D37").Value = '6.97'
$ws.Range("E37").Value = '  -2.42%  '

# Row 38
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").Value = '31.23'
$ws.Range("E38").Value = '  -2.36%  '

# Row 39
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '168.64'
$ws.Range("E39").Value = '  +0.29%  '

# Row 40
$ws.Range("D40").Value = '3.458.23'
$ws.Range("E40").Value = '  -1.12%  '

# Row 41
$ws.Range("E41").Value = '  +1.76%  '

# Row 42
$ws.Range("D42").Value = '42.67'
$ws.Range("E42").Value = '  +0.65%  '

# Row 43
$ws.Range("E43").Value = '  -1.39%  '

# Row 44
$ws.Range("D44").Value = '4.44'
$ws.Range("E44").Value = '  +1.73%  '

# Row 45
$ws.Range("E45").Value = '  -3.10%  '

# Row 46
$ws.Range("E46").Value = '  -3.78%  '

# Row 47
$ws.Range("D47").Value = '2.550.57'
$ws.Range("E47").Value = '  -1.58%  '

# Row 48
$ws.Range("E48").Value = '  +0.26%  '

# Row 49
$ws.Range("D49").Value = '22.55'
$ws.Range("E49").Value = '  -1.88%  '

# Row 51
$ws.Range("E51").Value = '  -6.02%  '
